# Auto-generated Excel COM-interop script to apply crypto price/volume updates
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, [string]$value)
    # Preserve the cell's existing style/format while forcing the incoming
    # text (even if it looks numeric, e.g. '0.994') to be stored as text,
    # matching the workbook's inlineStr cells instead of being parsed as a number.
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '60.922.27'
Set-TextValue $ws.Range('E2') '  -1.91%  '
Set-TextValue $ws.Range('D3') '2.417.49'
Set-TextValue $ws.Range('E3') '  -1.14%  '
Set-TextValue $ws.Range('D4') '0.994'
Set-TextValue $ws.Range('E4') '  -0.43%  '
Set-TextValue $ws.Range('D5') '570.19'
Set-TextValue $ws.Range('E5') '  -2.42%  '
Set-TextValue $ws.Range('D6') '139.66'
Set-TextValue $ws.Range('E6') '  -2.21%  '
Set-TextValue $ws.Range('E7') '  +0.12%  '
Set-TextValue $ws.Range('E8') '  -0.65%  '
Set-TextValue $ws.Range('D9') '2.401.09'
Set-TextValue $ws.Range('E9') '  -1.57%  '
Set-TextValue $ws.Range('E10') '  -1.51%  '
Set-TextValue $ws.Range('E11') '  -0.58%  '
Set-TextValue $ws.Range('D12') '5.07'
Set-TextValue $ws.Range('E12') '  -2.66%  '
Set-TextValue $ws.Range('D13') '0.338'
Set-TextValue $ws.Range('E13') '  -1.68%  '
Set-TextValue $ws.Range('D14') '26.15'
Set-TextValue $ws.Range('E14') '  -1.42%  '
Set-TextValue $ws.Range('E15') '  -2.90%  '
Set-TextValue $ws.Range('E16') '  -2.01%  '
Set-TextValue $ws.Range('D17') '60.702.74'
Set-TextValue $ws.Range('E17') '  -2.31%  '
Set-TextValue $ws.Range('D18') '2.392.82'
Set-TextValue $ws.Range('E18') '  -1.94%  '
Set-TextValue $ws.Range('D19') '7.78'
Set-TextValue $ws.Range('E19') '  +8.97%  '
Set-TextValue $ws.Range('D20') '10.65'
Set-TextValue $ws.Range('E20') '  -1.02%  '
Set-TextValue $ws.Range('D21') '323.01'
Set-TextValue $ws.Range('E21') '  -1.00%  '
Set-TextValue $ws.Range('E22') '  -1.56%  '
Set-TextValue $ws.Range('D23') '6.11'
Set-TextValue $ws.Range('E23') '  +2.20%  '
Set-TextValue $ws.Range('E24') '  +0.12%  '
Set-TextValue $ws.Range('D25') '1.84'
Set-TextValue $ws.Range('E25') '  -3.89%  '
Set-TextValue $ws.Range('D26') '64.84'
Set-TextValue $ws.Range('E26') '  -1.38%  '
Set-TextValue $ws.Range('D27') '581.53'
Set-TextValue $ws.Range('E27') '  -2.25%  '
Set-TextValue $ws.Range('E28') '  -9.45%  '
Set-TextValue $ws.Range('D29') '2.536.76'
Set-TextValue $ws.Range('E29') '  -1.17%  '
Set-TextValue $ws.Range('D30') '0.0₃0929'
Set-TextValue $ws.Range('E30') '  -4.47%  '
Set-TextValue $ws.Range('D31') '7.88'
Set-TextValue $ws.Range('E31') '  -1.28%  '
Set-TextValue $ws.Range('D32') '1.34'
Set-TextValue $ws.Range('E32') '  -5.23%  '
Set-TextValue $ws.Range('D33') '1.83'
Set-TextValue $ws.Range('E33') '  -3.08%  '
Set-TextValue $ws.Range('D34') '0.131'
Set-TextValue $ws.Range('E34') '  -3.23%  '
Set-TextValue $ws.Range('E35') '  -0.02%  '
Set-TextValue $ws.Range('D36') '1.40'
Set-TextValue $ws.Range('E36') '  -1.67%  '
Set-TextValue $ws.Range('D37') '4.60'
Set-TextValue $ws.Range('E37') '  -5.95%  '
Set-TextValue $ws.Range('D38') '151.24'
Set-TextValue $ws.Range('E38') '  -2.10%  '
Set-TextValue $ws.Range('E39') '  -1.94%  '
Set-TextValue $ws.Range('D40') '18.24'
Set-TextValue $ws.Range('E40') '  -0.84%  '
Set-TextValue $ws.Range('D41') '5.13'
Set-TextValue $ws.Range('E41') '  -2.98%  '
Set-TextValue $ws.Range('E42') '  +0.03%  '
Set-TextValue $ws.Range('E43') '  -2.24%  '
Set-TextValue $ws.Range('D44') '41.13'
Set-TextValue $ws.Range('E44') '  -5.03%  '
Set-TextValue $ws.Range('D45') '2.35'
Set-TextValue $ws.Range('E45') '  -6.68%  '
Set-TextValue $ws.Range('B46') 'Aave'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D46') '142.42'
Set-TextValue $ws.Range('E46') '  +0.40%  '
Set-TextValue $ws.Range('B47') 'BabyDogeCoin'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
Set-TextValue $ws.Range('D47') '0.0₆0267'
Set-TextValue $ws.Range('E47') '  +0.49%  '
Set-TextValue $ws.Range('D48') '3.50'
Set-TextValue $ws.Range('E48') '  -3.57%  '
Set-TextValue $ws.Range('E49') '  -2.39%  '
Set-TextValue $ws.Range('D50') '19.38'
Set-TextValue $ws.Range('E50') '  -2.53%  '
Set-TextValue $ws.Range('D51') '0.0502'
Set-TextValue $ws.Range('E51') '  -3.52%  '
